{"js": "// 1) Append a new run with the text \"\u0430\u0435\u0440\u043a\u0430\u0435\u0440\u0430\u0435\u0440\u043a\u0430\u0435\u0440\" right after the\n//    existing \"\u041c\u043e\u0441\u043a\u0432\u0430 202\" + \"3\" runs on the cover page, keeping the same\n//    run formatting (Times New Roman, 14pt / half-points 28, incl. the\n//    complex-script font & size so w:cs / w:szCs are emitted).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet targetParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"\u041c\u043e\u0441\u043a\u0432\u0430 2023\") {\n    targetParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (targetParagraph) {\n  const endRange = targetParagraph.getRange(Word.RangeLocation.end);\n  const ooxml =\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body><w:p><w:r><w:rPr>' +\n    '<w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/>' +\n    '<w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/>' +\n    '</w:rPr><w:t>\u0430\u0435\u0440\u043a\u0430\u0435\u0440\u0430\u0435\u0440\u043a\u0430\u0435\u0440</w:t></w:r></w:p></w:body></w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>';\n  endRange.insertOoxml(ooxml, Word.InsertLocation.end);\n  await context.sync();\n}\n\n// 2) Every run that directly hosts an inline picture gets marked\n//    \"do not spell/grammar-check\" (w:noProof) \u2014 this is what Word stamps\n//    on picture runs once the document has been interacted with.\nconst pictures = body.inlinePictures;\npictures.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < pictures.items.length; i++) {\n  const pictureRange = pictures.items[i].getRange();\n  pictureRange.hasNoProofing = true;\n}\nawait context.sync();\n", "ps1": "# 1) Append a new run with the text \"\u0430\u0435\u0440\u043a\u0430\u0435\u0440\u0430\u0435\u0440\u043a\u0430\u0435\u0440\" right after the\n#    existing \"\u041c\u043e\u0441\u043a\u0432\u0430 202\" + \"3\" runs on the cover page, keeping the same\n#    run formatting (Times New Roman, 14pt, incl. the complex-script\n#    font/size so w:cs / w:szCs get emitted on the new run as well).\n$d = $word.ActiveDocument\n\n$targetParagraph = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -eq \"\u041c\u043e\u0441\u043a\u0432\u0430 2023`r\") {\n        $targetParagraph = $p\n        break\n    }\n}\n\nif ($targetParagraph -ne $null) {\n    $r = $targetParagraph.Range\n    $r.InsertAfter(\"\u0430\u0435\u0440\u043a\u0430\u0435\u0440\u0430\u0435\u0440\u043a\u0430\u0435\u0440\")\n    $r.Font.Name = \"Times New Roman\"\n    $r.Font.Size = 14\n    $r.Font.NameBi = \"Times New Roman\"\n    $r.Font.SizeBi = 14\n}\n\n# 2) Every run that directly hosts an inline picture gets marked\n#    \"do not spell/grammar-check\" (w:noProof) - this is what Word stamps\n#    on picture runs once the document has been interacted with.\nforeach ($shp in $d.InlineShapes) {\n    $shp.Range.NoProofing = $true\n}\n"}
